$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "43.945.14"
Set-TextValue "E2" "  +3.07%  "
Set-TextValue "D3" "2.273.49"
Set-TextValue "E3" "  +2.91%  "
Set-TextValue "E4" "  -0.34%  "
Set-TextValue "D5" "321.26"
Set-TextValue "E5" "  +2.24%  "
Set-TextValue "D6" "103.10"
Set-TextValue "E6" "  +4.54%  "
Set-TextValue "D7" "0.588"
Set-TextValue "E7" "  +2.97%  "
Set-TextValue "E8" "  -0.18%  "
Set-TextValue "D9" "0.573"
Set-TextValue "E9" "  +3.02%  "
Set-TextValue "D10" "38.06"
Set-TextValue "E10" "  +4.31%  "
Set-TextValue "D11" "0.0843"
Set-TextValue "E11" "  +2.92%  "
Set-TextValue "D12" "7.87"
Set-TextValue "E12" "  +3.43%  "
Set-TextValue "E13" "  +2.96%  "
Set-TextValue "D14" "2.616.63"
Set-TextValue "E14" "  +2.93%  "
Set-TextValue "D15" "0.878"
Set-TextValue "E15" "  +3.67%  "
Set-TextValue "D16" "14.59"
Set-TextValue "E16" "  +4.73%  "
Set-TextValue "D17" "2.275.41"
Set-TextValue "E17" "  +3.87%  "
Set-TextValue "D18" "43.884.32"
Set-TextValue "E18" "  +3.22%  "
Set-TextValue "D19" "14.31"
Set-TextValue "E19" "  +3.27%  "
Set-TextValue "D20" "0.0₃0994"
Set-TextValue "E20" "  +4.34%  "
Set-TextValue "D21" "6.69"
Set-TextValue "E21" "  +3.59%  "
Set-TextValue "D22" "66.25"
Set-TextValue "E22" "  +1.41%  "
Set-TextValue "D23" "3.21"
Set-TextValue "E23" "  +0.67%  "
Set-TextValue "D24" "239.49"
Set-TextValue "E24" "  +2.69%  "
Set-TextValue "D25" "2.24"
Set-TextValue "E25" "  +6.03%  "
Set-TextValue "E26" "  +0.07%  "
Set-TextValue "E27" "  +3.39%  "
Set-TextValue "D28" "10.24"
Set-TextValue "E28" "  +1.72%  "
Set-TextValue "D29" "39.16"
Set-TextValue "E29" "  +16.81%  "
Set-TextValue "D30" "2.21"
Set-TextValue "E30" "  +2.12%  "
Set-TextValue "D31" "6.52"
Set-TextValue "E31" "  +2.20%  "
Set-TextValue "D32" "0.0887"
Set-TextValue "E32" "  +1.38%  "
Set-TextValue "D33" "20.57"
Set-TextValue "E33" "  +0.87%  "
Set-TextValue "D34" "160.70"
Set-TextValue "E34" "  +1.69%  "
Set-TextValue "B35" "LidoDAOToken"
Set-TextValue "C35" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D35" "3.41"
Set-TextValue "E35" "  +6.30%  "
Set-TextValue "B36" "WEMIXToken"
Set-TextValue "C36" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D36" "2.73"
Set-TextValue "E36" "  +1.15%  "
Set-TextValue "D37" "2.03"
Set-TextValue "E37" "  +9.63%  "
Set-TextValue "D38" "0.122"
Set-TextValue "E38" "  +0.70%  "
Set-TextValue "D39" "4.52"
Set-TextValue "E39" "  +3.07%  "
Set-TextValue "D40" "0.106"
Set-TextValue "E40" "  +4.10%  "
Set-TextValue "E41" "  +11.29%  "
Set-TextValue "D42" "15.64"
Set-TextValue "E42" "  +32.09%  "
Set-TextValue "D43" "0.0329"
Set-TextValue "E43" "  +3.39%  "
Set-TextValue "E44" "  -0.13%  "
Set-TextValue "D45" "1.822.09"
Set-TextValue "E45" "  +2.77%  "
Set-TextValue "D46" "0.209"
Set-TextValue "E46" "  +1.87%  "
Set-TextValue "D47" "86.60"
Set-TextValue "E47" "  -2.29%  "
Set-TextValue "D48" "5.42"
Set-TextValue "E48" "  +2.55%  "
Set-TextValue "D49" "76.86"
Set-TextValue "E49" "  +0.29%  "
Set-TextValue "D50" "8.92"
Set-TextValue "E50" "  +6.23%  "
Set-TextValue "D51" "59.95"
Set-TextValue "E51" "  +0.41%  "
